$d = $word.ActiveDocument

# 1) Add <w:noProof/> to the rPr of every run that contains a drawing (image).
#    These are paragraphs 5, 7, 13, 15, 16, 19, 21 (1-based, in the ORIGINAL
#    paragraph numbering, before any new paragraphs are inserted below).
$drawingParaIndexes = @(5, 7, 13, 15, 16, 19, 21)
foreach ($idx in $drawingParaIndexes) {
    $p = $d.Paragraphs($idx)
    $p.Range.NoProofing = 1
}

# 2) Move the <w:lastRenderedPageBreak/> marker from the last picture's run
#    (paragraph 21) to the "Output:" run (paragraph 18) - Word recalculated
#    pagination so the page break now renders before "Output:" instead of
#    before the final picture.
$outputPara = $d.Paragraphs(18)
$outputPara.Range.Find.Execute("Output:", $true, $false, $false, $false, $false, $true, 1, $false, "^l^&", 2)

